$d = $word.ActiveDocument

function Set-CellText($tableIndex, $row, $col, $newText) {
    $table = $d.Tables.Item($tableIndex)
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText
}

Set-CellText 1 1 4 "1.2.0 Dt : 07-Aug-2021"   # 1.1.0 Dt : 07-Aug-2021 -> 1.2.0 Dt : 07-Aug-2021
Set-CellText 1 2 4 "17-Aug-2021, 17:13"   # 07-Aug-2021, 15:16 -> 17-Aug-2021, 17:13
Set-CellText 2 2 4 "40"   # 35 -> 40
Set-CellText 3 2 3 "m³/hr"   # GPM -> m³/hr
Set-CellText 3 2 4 "24.1"   # 92.9 -> 24.1
Set-CellText 3 3 3 "°C"   # °F -> °C
Set-CellText 3 3 4 "12"   # 53.6 -> 12
Set-CellText 3 4 3 "°C"   # °F -> °C
Set-CellText 3 4 4 "7"   # 44.6 -> 7
Set-CellText 3 6 3 "mLC"   # ftLC -> mLC
Set-CellText 3 6 4 "5.3"   # 13.5 -> 5.3
Set-CellText 3 7 3 "DN"   # NPS -> DN
Set-CellText 3 7 4 "80"   # 3.2 -> 80
Set-CellText 3 10 3 "m² hr °C/kcal"   # ft² Hr °F/BTU -> m² hr °C/kcal
Set-CellText 3 11 3 "kg/cm²(g)"   # psi(g) -> kg/cm²(g)
Set-CellText 3 11 4 "8"   # 114 -> 8
Set-CellText 4 2 3 "kcal/Hr"   # MBH -> kcal/Hr
Set-CellText 4 2 4 "282403.8"   # 981.3 -> 282403.8
Set-CellText 4 3 3 "m³/hr"   # GPM -> m³/hr
Set-CellText 4 3 4 "48"   # 211.3 -> 48
Set-CellText 4 4 3 "°C"   # °F -> °C
Set-CellText 4 4 4 "32"   # 89.6 -> 32
Set-CellText 4 5 3 "°C"   # °F -> °C
Set-CellText 4 5 4 "37.9"   # 98.9 -> 37.9
Set-CellText 4 7 3 "m³/hr"   # GPM -> m³/hr
Set-CellText 4 8 3 "mLC"   # ftLC -> mLC
Set-CellText 4 8 4 "5"   # 16.4 -> 5
Set-CellText 4 9 3 "DN"   # NPS -> DN
Set-CellText 4 9 4 "100"   # 4 -> 100
Set-CellText 4 12 3 "m² hr °C/kcal"   # ft² Hr °F/BTU -> m² hr °C/kcal
Set-CellText 4 13 3 "kg/cm²(g)"   # psi(g) -> kg/cm²(g)
Set-CellText 4 13 4 "8"   # 114 -> 8
Set-CellText 5 2 3 "kcal/Hr"   # MBH -> kcal/Hr
Set-CellText 5 2 4 "161443.8"   # 561.3 -> 161443.8
Set-CellText 5 3 3 "m³/hr"   # GPM -> m³/hr
Set-CellText 5 3 4 "17.2"   # 43.4 -> 17.2
Set-CellText 5 4 3 "°C"   # °F -> °C
Set-CellText 5 4 4 "150"   # 239 -> 150
Set-CellText 5 5 3 "°C"   # °F -> °C
Set-CellText 5 5 4 "140"   # 212 -> 140
Set-CellText 5 6 4 "3"   # 4 -> 3
Set-CellText 5 7 3 "mLC"   # ftLC -> mLC
Set-CellText 5 7 4 "1.6"   # 12.1 -> 1.6
Set-CellText 5 8 3 "DN"   # NPS -> DN
Set-CellText 5 9 3 "kg/cm²(g)"   # psi(g) -> kg/cm²(g)
Set-CellText 5 9 4 "8"   # 113.8 -> 8
Set-CellText 7 2 3 "mm"   # in -> mm
Set-CellText 7 2 4 "2800"   # 111 -> 2800
Set-CellText 7 3 3 "mm"   # in -> mm
Set-CellText 7 3 4 "1450"   # 58 -> 1450
Set-CellText 7 4 3 "mm"   # in -> mm
Set-CellText 7 4 4 "2250"   # 89 -> 2250
Set-CellText 7 5 3 "ton"   # lbs -> ton
Set-CellText 7 5 4 "3.3"   # 7275.3 -> 3.3
Set-CellText 7 6 3 "ton"   # lbs -> ton
Set-CellText 7 6 4 "3.1"   # 6834.3 -> 3.1
Set-CellText 7 7 3 "ton"   # lbs -> ton
Set-CellText 7 7 4 "4.7"   # 10361.7 -> 4.7
Set-CellText 7 8 3 "ton"   # lbs -> ton
Set-CellText 7 8 4 "2.8"   # 6172.9 -> 2.8
Set-CellText 7 9 3 "mm"   # in -> mm
Set-CellText 7 9 4 "2700"   # 106.3 -> 2700

Write-Host "Done."